# Applies the "Tested changing null values" edit to the invoice export
# workbook: re-labels the Simple Fields / Items headers, fills in the new
# Billing/Shipping/Payment address + VAT fields, swaps several values
# around, clears a handful of now-empty fields, and tweaks a couple of
# item descriptions.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "Simple Fields"
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header row (row 1)
$ws1.Range("A1").Value2 = "Name"
$ws1.Range("B1").Value2 = "Vendor Address"
$ws1.Range("C1").Value2 = "Billing Name"
$ws1.Range("D1").Value2 = "Billing Address"
$ws1.Range("E1").Value2 = "Billing VAT Number"
$ws1.Range("F1").Value2 = "Shipping Address"
$ws1.Range("G1").Value2 = "Payment Address"
$ws1.Range("H1").Value2 = "Vendor VAT Number"
$ws1.Range("I1").Value2 = "Date"
$ws1.Range("J1").Value2 = "DueDate"
$ws1.Range("K1").Value2 = "Invoice Number"
$ws1.Range("L1").Value2 = "PO Number"
$ws1.Range("M1").Value2 = "Payment Terms"
$ws1.Range("N1").Value2 = "Shipping Charges"
$ws1.Range("O1").Value2 = "Tax Amount"
$ws1.Range("P1").Value2 = "Net Amount"
$ws1.Range("Q1").Value2 = "Total"
$ws1.Range("R1").Value2 = "Discount"
$ws1.Range("S1").Value2 = "Items"

# Value row (row 2)
$ws1.Range("A2").Value2 = "Sirius Cybernetics Corp."
$ws1.Range("B2").Value2 = "4592 Bell Street New York, NY 10018"
$ws1.Range("C2").Value2 = "CHOAM"
$ws1.Range("D2").Value2 = "27 Shield Wall Ave, Carthag, CH 1965 Arrakis"
$ws1.Range("E2").ClearContents()
$ws1.Range("F2").Value2 = "CHOAM 27 Shield Wall Ave, Carthag, CH 1965 Arrakis"
$ws1.Range("G2").Value2 = "Bill To"
$ws1.Range("H2").ClearContents()
$ws1.Range("I2").Value2 = "2016-01-27"
$ws1.Range("J2").Value2 = "2016-02-26"
$ws1.Range("K2").Value2 = "890127"
$ws1.Range("L2").ClearContents()
$ws1.Range("M2").Value2 = "30 days"
$ws1.Range("N2").ClearContents()
$ws1.Range("O2").Value2 = "47088.46"
$ws1.Range("P2").Value2 = "247834.00"
$ws1.Range("Q2").Value2 = "294922.46"
$ws1.Range("R2").ClearContents()
$ws1.Range("S2").Value2 = "table"

# ------------------------------------------------------------------
# Sheet 2: "Simple Fields - Formatted"
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Header row (row 1) - identical to sheet 1
$ws2.Range("A1").Value2 = "Name"
$ws2.Range("B1").Value2 = "Vendor Address"
$ws2.Range("C1").Value2 = "Billing Name"
$ws2.Range("D1").Value2 = "Billing Address"
$ws2.Range("E1").Value2 = "Billing VAT Number"
$ws2.Range("F1").Value2 = "Shipping Address"
$ws2.Range("G1").Value2 = "Payment Address"
$ws2.Range("H1").Value2 = "Vendor VAT Number"
$ws2.Range("I1").Value2 = "Date"
$ws2.Range("J1").Value2 = "DueDate"
$ws2.Range("K1").Value2 = "Invoice Number"
$ws2.Range("L1").Value2 = "PO Number"
$ws2.Range("M1").Value2 = "Payment Terms"
$ws2.Range("N1").Value2 = "Shipping Charges"
$ws2.Range("O1").Value2 = "Tax Amount"
$ws2.Range("P1").Value2 = "Net Amount"
$ws2.Range("Q1").Value2 = "Total"
$ws2.Range("R1").Value2 = "Discount"
$ws2.Range("S1").Value2 = "Items"

# Value row (row 2) - the address-style cells carry the CSV "Key,Value"
# breakdown and are wrapped (style s="1")
$ws2.Range("A2").Value2 = "Sirius Cybernetics Corp."

$ws2.Range("B2").Value2 = "Key,Value`n""Address Line 1"",""4592 Bell Street""`n""City"",""New York""`n""Country"",""United States""`n""State / County / Province"",""New York""`n""Zip Postal Code"",""10018"""
$ws2.Range("B2").WrapText = $true

$ws2.Range("C2").Value2 = "CHOAM"

$ws2.Range("D2").Value2 = "Key,Value`n""Address Line 1"",""27 Shield Ave Carthag CH Arrakis""`n""City"",""Wall""`n""Zip Postal Code"",""1965"""
$ws2.Range("D2").WrapText = $true

$ws2.Range("E2").ClearContents()

$ws2.Range("F2").Value2 = "Key,Value`n""Address Line 1"",""CHOAM 27 Shield Ave Carthag CH Arrakis""`n""City"",""Wall""`n""Zip Postal Code"",""1965"""
$ws2.Range("F2").WrapText = $true

$ws2.Range("G2").Value2 = "Key,Value`n""Address Line 1"",""To""`n""Address Line 2"",""""`n""Address Line 3"",""""`n""City"",""Bill""`n""Country"",""United States""`n""State / County / Province"",""Wyoming""`n""Zip Postal Code"",""""" 
$ws2.Range("G2").WrapText = $true

$ws2.Range("H2").ClearContents()
$ws2.Range("I2").Value2 = "2016-01-27"
$ws2.Range("J2").Value2 = "2016-02-26"
$ws2.Range("K2").Value2 = "890127"
$ws2.Range("L2").ClearContents()
$ws2.Range("M2").Value2 = "30 days"
$ws2.Range("N2").ClearContents()
$ws2.Range("O2").Value2 = "47088.46"
$ws2.Range("P2").Value2 = "247834.00"
$ws2.Range("Q2").Value2 = "294922.46"
$ws2.Range("R2").ClearContents()
$ws2.Range("S2").Value2 = "table"

# ------------------------------------------------------------------
# Sheet 3: "Items"
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Header row (row 1)
$ws3.Range("A1").Value2 = "Line Number"
$ws3.Range("B1").Value2 = "Description"
$ws3.Range("C1").Value2 = "Item PO Number"
$ws3.Range("D1").Value2 = "Quantity"
$ws3.Range("E1").Value2 = "Unit Price"
$ws3.Range("F1").Value2 = "Line Amount"
$ws3.Range("G1").Value2 = "Part Number"

# Item 1 (row 2)
$ws3.Range("A2").ClearContents()
$ws3.Range("B2").Value2 = "Nutrimatic Drinks Dispenser"
$ws3.Range("C2").ClearContents()
$ws3.Range("D2").Value2 = "10"
$ws3.Range("E2").Value2 = "4200.00"
$ws3.Range("F2").Value2 = "42000.00"

# Item 2 (row 3)
$ws3.Range("A3").ClearContents()
$ws3.Range("B3").ClearContents()
$ws3.Range("C3").Value2 = "Shipboard Computer ""Eddie"""
$ws3.Range("D3").Value2 = "17"
$ws3.Range("E3").Value2 = "8402.00"
$ws3.Range("F3").Value2 = "142834.00"

# Item 3 (row 4)
$ws3.Range("A4").ClearContents()
$ws3.Range("B4").Value2 = "Happy Vertical People Transporters"
$ws3.Range("C4").ClearContents()
$ws3.Range("D4").Value2 = "3"
$ws3.Range("E4").Value2 = "21000.00"
$ws3.Range("F4").Value2 = "63000.00"

# ------------------------------------------------------------------
# Sheet 4: "Items - Formatted"
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Header row (row 1) - identical to sheet 3
$ws4.Range("A1").Value2 = "Line Number"
$ws4.Range("B1").Value2 = "Description"
$ws4.Range("C1").Value2 = "Item PO Number"
$ws4.Range("D1").Value2 = "Quantity"
$ws4.Range("E1").Value2 = "Unit Price"
$ws4.Range("F1").Value2 = "Line Amount"
$ws4.Range("G1").Value2 = "Part Number"

# Item 1 (row 2)
$ws4.Range("A2").ClearContents()
$ws4.Range("B2").Value2 = "Nutrimatic Drinks Dispenser"
$ws4.Range("C2").ClearContents()
$ws4.Range("D2").Value2 = "10.00"
$ws4.Range("E2").Value2 = "4200.00"
$ws4.Range("F2").Value2 = "42000.00"

# Item 2 (row 3)
$ws4.Range("A3").ClearContents()
$ws4.Range("B3").ClearContents()
$ws4.Range("C3").Value2 = "Shipboard Computer ""Eddie"""
$ws4.Range("D3").Value2 = "17.00"
$ws4.Range("E3").Value2 = "8402.00"
$ws4.Range("F3").Value2 = "142834.00"

# Item 3 (row 4)
$ws4.Range("A4").ClearContents()
$ws4.Range("B4").Value2 = "Happy Vertical People Transporters"
$ws4.Range("C4").ClearContents()
$ws4.Range("D4").Value2 = "3.00"
$ws4.Range("E4").Value2 = "21000.00"
$ws4.Range("F4").Value2 = "63000.00"
